$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 79, shifting existing rows 79:109 down to 80:110
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new price-report record
$ws.Range("A79").Value = 9
$ws.Range("B79").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C79").Value = "Metropolitana"
$ws.Range("D79").Value = 45141
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = 100112035
$ws.Range("G79").Value = "Bruselas (repollito)"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 52
$ws.Range("K79").Value = 18000
$ws.Range("L79").Value = 18000
$ws.Range("M79").Value = 18000
$ws.Range("N79").Value = "$/malla 15 kilos"
$ws.Range("O79").Value = "Provincia de Quillota"
$ws.Range("P79").Value = 1200
$ws.Range("Q79").Value = 15
$ws.Range("R79").Value = "Hortaliza"
